# Menambahkan pengiriman notifikasi via telegram
#
# Inserts a new "Username Telegram" column between "Nama" and "Nomor",
# fills in the Telegram usernames, and reformats the phone numbers to
# the international +62 form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B -- shifts the old "Nomor" column (B) to C
# and the old "email" column (C) to D.
$ws.Columns.Item(2).Insert()

# New "Username Telegram" header + values.
$ws.Range("B1").Value = "Username Telegram"
$ws.Range("B2").Value = "hodsiador"
$ws.Range("B4").Value = "harugpa"

# Re-format the phone numbers with the +62 international prefix
# (column C keeps its original "store as text" formatting).
$ws.Range("C2").Value = "+6289907212341"
$ws.Range("C3").Value = "+6285156084242"
$ws.Range("C4").Value = "+6282456084203"

# The column insert does not re-target the existing hyperlinks, so
# rebuild them against their new home (column D).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:horange@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:woahe@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:haruna@gmail.com")
# Re-adding hyperlinks resets formatting -- restore the Hyperlink style.
$ws.Range("D2:D4").Style = "Hyperlink"

# Column widths (B/C/D).
$ws.Columns.Item(2).ColumnWidth = 18.666666666666668
$ws.Columns.Item(3).ColumnWidth = 14.833333333333334
$ws.Columns.Item(4).ColumnWidth = 18.5

# Match the author's final selection.
$ws.Range("B3").Select()
